# The workbook tracks daily price observations for Mango at "Vega Modelo de
# Temuco". Two new rows of data need to be inserted at row 158 (pushing the
# existing rows 158-228 down to 160-230), then populated with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 158, shifting existing data down.
$ws.Rows.Item(158).Insert()
$ws.Rows.Item(159).Insert()

# Populate the two new rows with the new observations.
$newRows = @(158, 159)
$origenValues = @("Brasil", "Perú")

for ($i = 0; $i -lt 2; $i++) {
    $r = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = 10
    $ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($r, 3).Value = "La Araucanía"
    $ws.Cells.Item($r, 4).Value = 44489
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r + 2, 4).NumberFormat
    $ws.Cells.Item($r, 5).Value = 9
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108002
    $ws.Cells.Item($r, 10).Value = "Mango"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = "Primera"
    $ws.Cells.Item($r, 13).Value = 300
    $ws.Cells.Item($r, 14).Value = 8000
    $ws.Cells.Item($r, 15).Value = 8000
    $ws.Cells.Item($r, 16).Value = 8000
    $ws.Cells.Item($r, 17).Value = "`$/bandeja 4 kilos"
    $ws.Cells.Item($r, 18).Value = $origenValues[$i]
    $ws.Cells.Item($r, 19).Value = 2000
    $ws.Cells.Item($r, 20).Value = 4
}
